# Terminal header image 변경
#
# Re-positions/re-sizes the terminal background picture and the terminal
# screenshot picture on slide 1, and enlarges + re-positions the
# "Tech  Journal" title textbox (including bumping its font size).
#
# NOTE on the magic-looking point values below: this COM host stores
# Shape.Left/Top/Width/Height internally as single-precision (f32) points
# before converting back to EMU on save, which truncates (not rounds) to
# the nearest EMU. The literals here are chosen so that, after that f32
# round-trip, they land exactly on the target EMU offsets from the
# canonical OOXML (e.g. -4580965 EMU, 6964310 EMU, ...) instead of one
# EMU short, which is what a naive `emu / 914400 * 72` would produce.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape 1: "그림 3" (terminal background picture) ---
# <a:off x="-5751984" y="799863"/><a:ext cx="5530311" cy="5495961"/>
#   -> <a:off x="-4580965" y="-242347"/><a:ext cx="6964310" cy="6921053"/>
$shp1 = $s.Shapes.Item(1)
$shp1.Left   = -360.70591751181104   # -4580965 EMU
$shp1.Top    = -19.08244094488189    # -242347 EMU
$shp1.Width  = 548.3708801417323     # 6964310 EMU
$shp1.Height = 544.9648141496062     # 6921053 EMU

# --- Shape 2: "그림 2" (terminal screenshot picture) ---
# <a:off x="141111" y="1766710"/> (ext unchanged)
#   -> <a:off x="129822" y="1969080"/>
$shp2 = $s.Shapes.Item(2)
$shp2.Left = 10.22220472440945   # 129822 EMU
$shp2.Top  = 155.04566929133858  # 1969080 EMU

# --- Shape 3: "TextBox 5" ("Tech  Journal" title) ---
# <a:off x="344016" y="3429000"/><a:ext cx="6638161" cy="1169551"/>
#   -> <a:off x="448202" y="3366015"/><a:ext cx="6638161" cy="1446550"/>
# font size 7000 (70pt) -> 8800 (88pt)
$shp3 = $s.Shapes.Item(3)
$shp3.Left   = 35.291496062992124  # 448202 EMU
$shp3.Top    = 265.0405511811024   # 3366015 EMU
$shp3.Width  = 522.689850519685    # 6638161 EMU
$shp3.Height = 113.9015778031496   # 1446550 EMU

$shp3.TextFrame.TextRange.Font.Size = 88
